$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> 1-based index map (A=1 ... V=22)
$colIndex = @{
    A=1; B=2; C=3; D=4; E=5; F=6; G=7; H=8; I=9; J=10; K=11; L=12;
    M=13; N=14; O=15; P=16; Q=17; R=18; S=19; T=20; U=21; V=22
}

# --- Part 1: rows whose match data (columns F:V) needs to be updated/re-shuffled ---
# (Index/country/tournament/season/date in columns A:E are unchanged for these rows.)
$affectedRows = @(
    @{ Row=4; F='AFC'; G=4; H='Jong Almere City'; I=1; J=1.35; K='19/08/2023 09:28'; L=1.36; M='19/08/2023 14:17'; N=5.49; O='19/08/2023 09:28'; P=5.43; Q='19/08/2023 14:17'; R=5.99; S='19/08/2023 09:28'; T=6; U='19/08/2023 14:17'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/afc-jong-almere-city/nLW387jH/' },
    @{ Row=5; F='Rijnsburgse Boys'; G=2; H='Kozakken Boys'; I=1; J=1.46; K='18/08/2023 03:13'; L=1.51; M='19/08/2023 00:50'; N=4.52; O='18/08/2023 03:13'; P=4.59; Q='19/08/2023 13:05'; R=4.5; S='18/08/2023 03:13'; T=4.75; U='19/08/2023 00:50'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/rijnsburgse-boys-kozakken-boys/hMweATL4/' },
    @{ Row=10; F='GVVV'; G=1; H='De Treffers'; I=1; J=2.11; K='25/08/2023 02:42'; L=2.4; M='26/08/2023 12:42'; N=3.44; O='25/08/2023 02:42'; P=3.78; Q='26/08/2023 12:42'; R=2.75; S='25/08/2023 02:42'; T=2.48; U='26/08/2023 12:42'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/gvvv-de-treffers/vikB2xKe/' },
    @{ Row=11; F='Noordwijk'; G=1; H='AFC'; I=1; J=3.27; K='25/08/2023 02:42'; L=2.9; M='26/08/2023 13:00'; N=3.67; O='25/08/2023 02:42'; P=3.6; Q='26/08/2023 13:00'; R=1.85; S='25/08/2023 02:42'; T=2.16; U='26/08/2023 13:00'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/vv-noordwijk-afc/jTXYLQKG/' },
    @{ Row=12; F='Lisse'; G=1; H='Excelsior Maassluis'; I=1; J=2.14; K='25/08/2023 02:42'; L=2.17; M='26/08/2023 14:23'; N=3.54; O='25/08/2023 02:42'; P=3.79; Q='26/08/2023 14:23'; R=2.72; S='25/08/2023 02:42'; T=2.78; U='26/08/2023 14:23'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/lisse-excelsior-maassluis/pG2YNnj4/' },
    @{ Row=13; F='Hardenberg'; G=3; H='ACV Assen'; I=0; J=1.64; K='26/08/2023 13:30'; L=1.6; M='26/08/2023 13:43'; N=5.05; O='26/08/2023 13:30'; P=5.1; Q='26/08/2023 14:29'; R=3.61; S='26/08/2023 13:30'; T=3.79; U='26/08/2023 14:29'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/hardenberg-acv-assen/6P3UOSzb/' },
    @{ Row=18; F='Excelsior Maassluis'; G=5; H='Jong Almere City'; I=3; J=1.83; K='02/09/2023 11:12'; L=1.74; M='02/09/2023 13:45'; N=3.99; O='02/09/2023 11:12'; P=3.99; Q='02/09/2023 13:49'; R=3.29; S='02/09/2023 11:12'; T=3.38; U='02/09/2023 13:42'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/excelsior-maassluis-jong-almere-city/hfFpHyZl/' },
    @{ Row=19; F='Hardenberg'; G=2; H='Jong Sparta Rotterdam'; I=0; J=1.43; K='01/09/2023 02:42'; L=1.43; M='02/09/2023 14:27'; N=4.67; O='01/09/2023 02:42'; P=5.34; Q='02/09/2023 14:27'; R=4.6; S='01/09/2023 02:42'; T=5.1; U='02/09/2023 14:27'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/hardenberg-jong-sparta-rotterdam/6cscQuC2/' },
    @{ Row=20; F='ACV Assen'; G=2; H='Lisse'; I=0; J=2.08; K='02/09/2023 11:12'; L=1.85; M='02/09/2023 14:21'; N=3.74; O='02/09/2023 11:12'; P=3.85; Q='02/09/2023 14:28'; R=2.83; S='02/09/2023 11:12'; T=3.49; U='02/09/2023 14:28'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/acv-assen-lisse/WpElGele/' },
    @{ Row=21; F='Scheveningen'; G=3; H='GVVV'; I=1; J=2.87; K='01/09/2023 02:42'; L=2.97; M='02/09/2023 12:56'; N=3.43; O='01/09/2023 02:42'; P=3.32; Q='02/09/2023 13:02'; R=2.09; S='01/09/2023 02:42'; T=2.24; U='02/09/2023 13:02'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/svv-scheveningen-gvvv/IslF1dZ1/' },
    @{ Row=25; F='Lisse'; G=0; H='Hardenberg'; I=5; J=4.45; K='08/09/2023 02:42'; L=4.19; M='09/09/2023 14:15'; N=4.34; O='08/09/2023 02:42'; P=5.2; Q='09/09/2023 14:15'; R=1.48; S='08/09/2023 02:42'; T=1.53; U='09/09/2023 14:15'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/lisse-hardenberg/fmRvRC2D/' },
    @{ Row=26; F='Noordwijk'; G=1; H='Excelsior Maassluis'; I=1; J=1.76; K='08/09/2023 02:42'; L=1.69; M='09/09/2023 14:27'; N=3.77; O='08/09/2023 02:42'; P=4.29; Q='09/09/2023 14:28'; R=3.37; S='08/09/2023 02:42'; T=3.81; U='09/09/2023 14:27'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/vv-noordwijk-excelsior-maassluis/6yYPlGfK/' },
    @{ Row=27; F='GVVV'; G=1; H='Rijnsburgse Boys'; I=4; J=2.22; K='08/09/2023 02:42'; L=2.43; M='09/09/2023 14:15'; N=3.5; O='08/09/2023 02:42'; P=3.49; Q='09/09/2023 14:15'; R=2.55; S='08/09/2023 02:42'; T=2.59; U='09/09/2023 14:15'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/gvvv-rijnsburgse-boys/bZC4sYnf/' },
    @{ Row=33; F='ACV Assen'; G=1; H='Noordwijk'; I=1; J=2.83; K='16/09/2023 12:42'; L=3.25; M='16/09/2023 14:28'; N=3.64; O='16/09/2023 12:42'; P=3.73; Q='16/09/2023 14:28'; R=2.19; S='16/09/2023 12:42'; T=1.96; U='16/09/2023 14:28'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/acv-assen-vv-noordwijk/E3lWM063/' },
    @{ Row=34; F='Excelsior Maassluis'; G=0; H='Kozakken Boys'; I=0; J=2.23; K='15/09/2023 02:42'; L=2.47; M='16/09/2023 14:26'; N=3.49; O='15/09/2023 02:42'; P=3.91; Q='16/09/2023 14:26'; R=2.56; S='15/09/2023 02:42'; T=2.36; U='16/09/2023 14:26'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/excelsior-maassluis-kozakken-boys/zDkSNtic/' },
    @{ Row=35; F='Hardenberg'; G=1; H='Jong Almere City'; I=0; J=1.29; K='16/09/2023 12:42'; L=1.3; M='16/09/2023 13:33'; N=6.35; O='16/09/2023 12:42'; P=6.28; Q='16/09/2023 13:33'; R=6.21; S='16/09/2023 12:42'; T=6.3; U='16/09/2023 13:33'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/hardenberg-jong-almere-city/00Od81iA/' },
    @{ Row=37; F='Lisse'; G=4; H='Jong Sparta Rotterdam'; I=7; J=2.56; K='15/09/2023 02:42'; L=3.28; M='16/09/2023 14:27'; N=3.53; O='15/09/2023 02:42'; P=3.53; Q='16/09/2023 14:26'; R=2.2; S='15/09/2023 02:42'; T=1.84; U='16/09/2023 14:27'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/lisse-jong-sparta-rotterdam/IyJ17L6G/' },
    @{ Row=59; F='Noordwijk'; G=0; H='Jong Almere City'; I=0; J=1.61; K='06/10/2023 01:42'; L=1.66; M='07/10/2023 13:44'; N=4.11; O='06/10/2023 01:42'; P=4.27; Q='07/10/2023 13:44'; R=3.79; S='06/10/2023 01:42'; T=4.03; U='07/10/2023 13:44'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/vv-noordwijk-jong-almere-city/lpshi9GI/' },
    @{ Row=60; F='GVVV'; G=2; H='Excelsior Maassluis'; I=1; J=1.43; K='06/10/2023 01:42'; L=1.62; M='06/10/2023 11:53'; N=4.65; O='06/10/2023 01:42'; P=4.46; Q='07/10/2023 12:33'; R=4.62; S='06/10/2023 01:42'; T=4.02; U='06/10/2023 19:05'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/gvvv-excelsior-maassluis/belufBWa/' },
    @{ Row=61; F='Spakenburg'; G=1; H='Hardenberg'; I=0; J=2.84; K='06/10/2023 02:12'; L=2.47; M='07/10/2023 14:45'; N=3.57; O='06/10/2023 02:12'; P=3.55; Q='07/10/2023 14:45'; R=2.02; S='06/10/2023 02:12'; T=2.51; U='07/10/2023 14:45'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/spakenburg-hardenberg/Cdc3bQNo/' },
    @{ Row=62; F='Quick Boys'; G=0; H='ACV Assen'; I=3; J=1.53; K='07/10/2023 11:42'; L=1.69; M='07/10/2023 14:56'; N=4.56; O='07/10/2023 11:42'; P=4.28; Q='07/10/2023 14:57'; R=4.34; S='07/10/2023 11:42'; T=3.85; U='07/10/2023 14:57'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/quick-boys-acv-assen/jmd7cpwh/' },
    @{ Row=63; F='ADO 20 Heemskerk'; G=2; H='Katwijk'; I=1; J=3.26; K='07/10/2023 11:42'; L=3.63; M='07/10/2023 14:46'; N=3.8; O='07/10/2023 11:42'; P=4.06; Q='07/10/2023 14:46'; R=1.88; S='07/10/2023 11:42'; T=1.78; U='07/10/2023 14:46'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/ado-20-heemskerk-katwijk/SUhqgVo6/' },
    @{ Row=64; F='Kozakken Boys'; G=2; H='Lisse'; I=0; J=1.42; K='06/10/2023 02:12'; L=1.51; M='07/10/2023 08:14'; N=4.45; O='06/10/2023 02:12'; P=4.41; Q='07/10/2023 13:02'; R=5.2; S='06/10/2023 02:12'; T=4.96; U='07/10/2023 12:44'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/kozakken-boys-lisse/bg9LZtoo/' },
    @{ Row=67; F='Noordwijk'; G=0; H='Jong Sparta Rotterdam'; I=2; J=1.94; K='13/10/2023 01:42'; L=1.95; M='14/10/2023 14:25'; N=3.69; O='13/10/2023 01:42'; P=4.09; Q='14/10/2023 14:25'; R=2.92; S='13/10/2023 01:42'; T=3.04; U='14/10/2023 14:25'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/vv-noordwijk-jong-sparta-rotterdam/GfsXit7T/' },
    @{ Row=68; F='Lisse'; G=1; H='Spakenburg'; I=3; J=3.93; K='13/10/2023 01:42'; L=5.65; M='14/10/2023 14:17'; N=4.04; O='13/10/2023 01:42'; P=4.6; Q='14/10/2023 14:17'; R=1.62; S='13/10/2023 01:42'; T=1.45; U='14/10/2023 14:17'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/lisse-spakenburg/fugOg2xH/' },
    @{ Row=69; F='Hardenberg'; G=1; H='Quick Boys'; I=3; J=2; K='13/10/2023 01:42'; L=2.18; M='14/10/2023 14:15'; N=3.66; O='13/10/2023 01:42'; P=4; Q='14/10/2023 14:16'; R=2.82; S='13/10/2023 01:42'; T=2.66; U='14/10/2023 14:14'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/hardenberg-quick-boys/txkKfrNA/' },
    @{ Row=70; F='ACV Assen'; G=3; H='GVVV'; I=0; J=2.7; K='14/10/2023 12:34'; L=2.77; M='14/10/2023 14:17'; N=3.82; O='14/10/2023 12:34'; P=3.49; Q='14/10/2023 14:17'; R=2.21; S='14/10/2023 12:34'; T=2.29; U='14/10/2023 14:17'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/acv-assen-gvvv/M9Ph9sx4/' },
    @{ Row=76; F='Scheveningen'; G=0; H='Katwijk'; I=4; J=3.32; K='27/10/2023 02:42'; L=2.94; M='28/10/2023 14:26'; N=3.55; O='27/10/2023 02:42'; P=3.74; Q='28/10/2023 14:26'; R=1.83; S='27/10/2023 02:42'; T=2.1; U='28/10/2023 14:26'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/svv-scheveningen-katwijk/dSwyMKL9/' },
    @{ Row=77; F='GVVV'; G=2; H='Hardenberg'; I=1; J=2.57; K='27/10/2023 02:42'; L=2.53; M='28/10/2023 14:19'; N=3.4; O='27/10/2023 02:42'; P=3.58; Q='28/10/2023 14:19'; R=2.25; S='27/10/2023 02:42'; T=2.44; U='28/10/2023 14:19'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/gvvv-hardenberg/Q3qCDdTk/' },
    @{ Row=78; F='Spakenburg'; G=4; H='Jong Almere City'; I=3; J=1.46; K='27/10/2023 03:12'; L=1.43; M='28/10/2023 09:52'; N=4.52; O='27/10/2023 03:12'; P=4.97; Q='28/10/2023 13:03'; R=4.5; S='27/10/2023 03:12'; T=5.33; U='28/10/2023 09:52'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/spakenburg-jong-almere-city/rqtmJI5S/' },
    @{ Row=80; F='Kozakken Boys'; G=1; H='Noordwijk'; I=3; J=2.26; K='27/10/2023 03:12'; L=2.49; M='28/10/2023 14:58'; N=3.34; O='27/10/2023 03:12'; P=3.47; Q='28/10/2023 14:58'; R=2.6; S='27/10/2023 03:12'; T=2.53; U='28/10/2023 14:58'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/kozakken-boys-vv-noordwijk/4hsqKbjM/' }
)

foreach ($item in $affectedRows) {
    $r = $item.Row
    foreach ($col in @('F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T','U','V')) {
        $ws.Cells.Item($r, $colIndex[$col]).Value = $item[$col]
    }
}

# --- Part 2: brand new rows 83:89 appended at the bottom of the sheet ---
$newRows = @(
    @{ Row=83; A=82; B='netherlands'; C='tweede-divisie'; D='2023-2024'; E=45234.60416666666; F='ACV Assen'; G=2; H='De Treffers'; I=3; J=3.41; K='03/11/2023 02:42'; L=3.57; M='04/11/2023 14:09'; N=3.78; O='03/11/2023 02:42'; P=3.88; Q='04/11/2023 14:09'; R=1.78; S='03/11/2023 02:42'; T=1.83; U='04/11/2023 14:09'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/acv-assen-de-treffers/dQ2a8M87/' },
    @{ Row=84; A=83; B='netherlands'; C='tweede-divisie'; D='2023-2024'; E=45234.60416666666; F='Excelsior Maassluis'; G=2; H='Scheveningen'; I=2; J=2.84; K='03/11/2023 02:42'; L=3.42; M='04/11/2023 14:26'; N=3.57; O='03/11/2023 02:42'; P=3.65; Q='04/11/2023 14:26'; R=2.02; S='03/11/2023 02:42'; T=1.93; U='04/11/2023 14:26'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/excelsior-maassluis-svv-scheveningen/IZ3e92g1/' },
    @{ Row=85; A=84; B='netherlands'; C='tweede-divisie'; D='2023-2024'; E=45234.60416666666; F='Lisse'; G=2; H='GVVV'; I=2; J=3.81; K='03/11/2023 02:42'; L=4.56; M='04/11/2023 14:14'; N=4; O='03/11/2023 02:42'; P=4.35; Q='04/11/2023 14:14'; R=1.62; S='03/11/2023 02:42'; T=1.57; U='04/11/2023 14:14'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/lisse-gvvv/pSG8I6WE/' },
    @{ Row=86; A=85; B='netherlands'; C='tweede-divisie'; D='2023-2024'; E=45234.60416666666; F='Noordwijk'; G=0; H='Spakenburg'; I=3; J=2.42; K='03/11/2023 02:42'; L=2.46; M='04/11/2023 14:22'; N=3.44; O='03/11/2023 02:42'; P=3.62; Q='04/11/2023 14:22'; R=2.42; S='03/11/2023 02:42'; T=2.49; U='04/11/2023 14:22'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/vv-noordwijk-spakenburg/6s5HGp1R/' },
    @{ Row=87; A=86; B='netherlands'; C='tweede-divisie'; D='2023-2024'; E=45234.625; F='Kozakken Boys'; G=1; H='Jong Sparta Rotterdam'; I=2; J=2.16; K='03/11/2023 03:12'; L=2.77; M='04/11/2023 14:57'; N=3.57; O='03/11/2023 03:12'; P=3.74; Q='04/11/2023 14:57'; R=2.61; S='03/11/2023 03:12'; T=2.19; U='04/11/2023 14:57'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/kozakken-boys-jong-sparta-rotterdam/SpBrC49r/' },
    @{ Row=88; A=87; B='netherlands'; C='tweede-divisie'; D='2023-2024'; E=45234.64583333334; F='Jong Almere City'; G=2; H='Quick Boys'; I=0; J=3.94; K='03/11/2023 03:43'; L=4.25; M='04/11/2023 15:27'; N=4.3; O='03/11/2023 03:43'; P=4.8; Q='04/11/2023 15:27'; R=1.56; S='03/11/2023 03:43'; T=1.56; U='04/11/2023 15:27'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/jong-almere-city-quick-boys/OIFCHQoL/' },
    @{ Row=89; A=88; B='netherlands'; C='tweede-divisie'; D='2023-2024'; E=45234.64583333334; F='Katwijk'; G=3; H='Rijnsburgse Boys'; I=0; J=1.65; K='03/11/2023 03:43'; L=1.47; M='04/11/2023 15:21'; N=3.92; O='03/11/2023 03:43'; P=4.64; Q='04/11/2023 15:21'; R=3.73; S='03/11/2023 03:43'; T=5.22; U='04/11/2023 15:21'; V='https://www.betexplorer.com/football/netherlands/tweede-divisie/katwijk-rijnsburgse-boys/0b9jArve/' }
)

foreach ($item in $newRows) {
    $r = $item.Row
    # Copy formatting (styles) from the last pre-existing data row (82) down to the new row
    $ws.Range("A82:V82").Copy()
    $ws.Range("A" + $r + ":V" + $r).PasteSpecial(-4122)

    foreach ($col in @('A','B','C','D','E','F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T','U','V')) {
        $ws.Cells.Item($r, $colIndex[$col]).Value = $item[$col]
    }
}

Write-Host "Edit complete"
